$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = "Jacopo Chemini"
$ws.Range("B28").Value = "Stefano  Tita | Clitoriders"
$ws.Range("C28").Value = "Daniel Pedrotti | IMONTAGNA"
$ws.Range("D28").Value = "Michele Merighi | Clitoriders"
$ws.Range("E28").Value = "maikol  azocar | Mai una gioia"
$ws.Range("F28").Value = "Emanuele Toss | 4SINS"
